# Update column F ("dSF") values on Sheet1 per the repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    3  = -1
    4  = 3
    5  = 6
    6  = -5
    7  = -3
    8  = -3
    9  = -6
    10 = 4
    11 = 1
    12 = -4
    13 = -7
    14 = 1
    15 = -3
    18 = -5
    20 = 1
    21 = -2
    22 = -6
    23 = 1
    24 = -2
    25 = -1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
